# Apply updated cryptos list values (Price + Volume(1h) columns).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" values are plain decimal-looking strings (e.g. "330.72") that
# Excel would otherwise auto-convert to numbers on assignment. Force those
# specific cells to a Text number format first so the literal string is kept,
# matching the inline string cells already in the workbook.
$textForceCells = @('D4', 'D5', 'D7', 'D8', 'D9', 'D10', 'D11', 'D12', 'D14', 'D16', 'D17', 'D18', 'D19', 'D20', 'D21', 'D22', 'D24', 'D25', 'D26', 'D27', 'D28', 'D29', 'D30', 'D32', 'D33', 'D34', 'D35', 'D37', 'D38', 'D39', 'D40', 'D41', 'D42', 'D44', 'D45', 'D46', 'D47', 'D48', 'D49', 'D50', 'D51')
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '29.599.24'
$ws.Range('E2').Value = '  -2.56%  '
$ws.Range('D3').Value = '2.000.30'
$ws.Range('E3').Value = '  -5.03%  '
$ws.Range('D4').Value = '1.015'
$ws.Range('E4').Value = '  +0.75%  '
$ws.Range('D5').Value = '330.72'
$ws.Range('E5').Value = '  -4.12%  '
$ws.Range('E6').Value = '  +0.64%  '
$ws.Range('D7').Value = '0.4996'
$ws.Range('E7').Value = '  -4.46%  '
$ws.Range('D8').Value = '0.4242'
$ws.Range('E8').Value = '  -4.38%  '
$ws.Range('D9').Value = '54.65'
$ws.Range('E9').Value = '  -0.06%  '
$ws.Range('D10').Value = '0.09008'
$ws.Range('E10').Value = '  -3.79%  '
$ws.Range('D11').Value = '1.120'
$ws.Range('E11').Value = '  -4.39%  '
$ws.Range('D12').Value = '23.40'
$ws.Range('E12').Value = '  -6.03%  '
$ws.Range('D13').Value = '2.007.75'
$ws.Range('E13').Value = '  -4.64%  '
$ws.Range('D14').Value = '8.079'
$ws.Range('E14').Value = '  -6.69%  '
$ws.Range('E15').Value = '  -6.08%  '
$ws.Range('D16').Value = '1.015'
$ws.Range('E16').Value = '  +0.65%  '
$ws.Range('D17').Value = '94.51'
$ws.Range('E17').Value = '  -7.21%  '
$ws.Range('D18').Value = '0.00001112'
$ws.Range('E18').Value = '  -4.21%  '
$ws.Range('D19').Value = '0.06672'
$ws.Range('E19').Value = '  -0.70%  '
$ws.Range('D20').Value = '19.80'
$ws.Range('E20').Value = '  -6.64%  '
$ws.Range('D21').Value = '1.013'
$ws.Range('E21').Value = '  +0.75%  '
$ws.Range('D22').Value = '5.955'
$ws.Range('E22').Value = '  -6.29%  '
$ws.Range('D23').Value = '29.624.24'
$ws.Range('E23').Value = '  -2.64%  '
$ws.Range('D24').Value = '12.02'
$ws.Range('E24').Value = '  -4.68%  '
$ws.Range('D25').Value = '2.282'
$ws.Range('E25').Value = '  -0.77%  '
$ws.Range('D26').Value = '158.74'
$ws.Range('E26').Value = '  -2.42%  '
$ws.Range('D27').Value = '20.71'
$ws.Range('E27').Value = '  -5.70%  '
$ws.Range('D28').Value = '6.382'
$ws.Range('E28').Value = '  -5.59%  '
$ws.Range('D29').Value = '2.307'
$ws.Range('E29').Value = '  -8.39%  '
$ws.Range('D30').Value = '128.55'
$ws.Range('E30').Value = '  -3.87%  '
$ws.Range('E31').Value = '  -7.54%  '
$ws.Range('D32').Value = '0.09938'
$ws.Range('E32').Value = '  -5.61%  '
$ws.Range('D33').Value = '1.569'
$ws.Range('E33').Value = '  -7.17%  '
$ws.Range('D34').Value = '5.845'
$ws.Range('E34').Value = '  -6.45%  '
$ws.Range('D35').Value = '3.814'
$ws.Range('E35').Value = '  -2.81%  '
$ws.Range('E36').Value = '  -8.04%  '
$ws.Range('D37').Value = '0.02469'
$ws.Range('E37').Value = '  -6.31%  '
$ws.Range('D38').Value = '1.307'
$ws.Range('E38').Value = '  -3.26%  '
$ws.Range('D39').Value = '0.06359'
$ws.Range('E39').Value = '  -6.24%  '
$ws.Range('D40').Value = '0.6565'
$ws.Range('E40').Value = '  -6.65%  '
$ws.Range('D41').Value = '11.70'
$ws.Range('E41').Value = '  -6.74%  '
$ws.Range('D42').Value = '0.2058'
$ws.Range('E42').Value = '  -7.47%  '
$ws.Range('E43').Value = '  +0.65%  '
$ws.Range('D44').Value = '0.6345'
$ws.Range('E44').Value = '  -7.36%  '
$ws.Range('D45').Value = '13.52'
$ws.Range('E45').Value = '  -6.80%  '
$ws.Range('D46').Value = '2.207'
$ws.Range('E46').Value = '  -6.21%  '
$ws.Range('D47').Value = '1.291'
$ws.Range('E47').Value = '  -7.71%  '
$ws.Range('D48').Value = '3.523'
$ws.Range('E48').Value = '  -3.30%  '
$ws.Range('D49').Value = '0.00000000339'
$ws.Range('E49').Value = '  -3.02%  '
$ws.Range('D50').Value = '0.06988'
$ws.Range('E50').Value = '  -3.60%  '
$ws.Range('D51').Value = '1.125'
$ws.Range('E51').Value = '  -7.20%  '
